# Create the "UserStories" worksheet after Sheet1
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "UserStories"

# Header row
$headers = @("Identifier", "Priority", "Overview")
$col = 2
foreach ($h in $headers) {
    $ws2.Cells.Item(2, $col).Value = $h
    $col = $col + 1
}

# Data rows: Identifier, Priority, Overview
$data = @(
    @("STORY-01", "High",   "As someone that takes on a lot of self-directed work, I'd like to plan & record sessions of work."),
    @("STORY-02", "Medium", "As a freelancer I want to track how much time I spend on a given project."),
    @("STORY-03", "Low",    "As a student I want to keep track of work on different modules at the same time."),
    @("STORY-04", "High",   "As someone that gets distracted easily, I want to set myself short-term objectives as an incentive to tackle my work."),
    @("STORY-05", "Low",    "As a freelancer I want access to a backlog of tasks I checked off on a given project."),
    @("STORY-06", "Low",    "As someone with data analysis skills, I want to be able to export data from the system in a .csv format"),
    @("STORY-07", "Low",    "As someone with many ongoing projects to keep track of I want a hierarchical tagging/categories system"),
    @("STORY-08", "Low",    "As someone self-employed I want to search back through my sessions of work to find when I was working on certain tasks."),
    @("STORY-09", "Low",    "As a student that uses colour-encoding to organise my work I want to be able to customise colour of some ui elements/categories.")
)

$row = 3
foreach ($rec in $data) {
    $ws2.Cells.Item($row, 2).Value = $rec[0]
    $ws2.Cells.Item($row, 3).Value = $rec[1]
    $ws2.Cells.Item($row, 4).Value = $rec[2]
    $row = $row + 1
}

# Formatting: columns B:D rows 2-11, centered + thin border
$range = $ws2.Range("B2:D11")
$range.HorizontalAlignment = -4108
$range.Borders.LineStyle = 1

# Overview column data cells (D3:D11) use the "Good" (green) cell style
$good = $ws2.Range("D3:D11")
$good.Style = "Good"
$good.Borders.LineStyle = 1
$good.HorizontalAlignment = -4108

# Column widths
$ws2.Columns.Item("B").ColumnWidth = 15.42578125
$ws2.Columns.Item("C").ColumnWidth = 15.42578125
$ws2.Columns.Item("D").ColumnWidth = 122.140625
$ws2.Columns.Item("E").ColumnWidth = 12.140625

# Selection / view state
$ws2.Range("D15").Select()
$ws1.Range("L23").Select()

$wb.Worksheets.Item(1).Activate()
$ws2.Activate()
